$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6075708866119385
$ws.Range("B1").Value = 2.102619886398315
$ws.Range("D1").Value = 1.160134673118591
$ws.Range("E1").Value = 1.248434901237488
